$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the input values (column B) ---
$ws.Range("B1").Value = 60
$ws.Range("B2").Value = 200
$ws.Range("B3").Value = 7
$ws.Range("B4").Value = 2.5
$ws.Range("B5").Value = 2

# --- Apply a custom "0.0" number format ---
# D (calculated/output) column first, then B (input) column -- this order
# matches how the two new cell formats were introduced in the workbook.
$ws.Range("D1:D5").NumberFormat = "0.0"
$ws.Range("B1:B5").NumberFormat = "0.0"

# --- Resize columns B and D slightly wider ---
$ws.Columns.Item(2).ColumnWidth = 5.498697916666667
$ws.Columns.Item(4).ColumnWidth = 5.166666666666667

# --- Move the selection ---
$ws.Range("C11").Select()
